$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.042.10'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '3.152.06'
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('E4').Value = '  +0.18%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.11'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +0.48%  '
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.21'
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.143.33'
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('E10').Value = '  +0.78%  '
$style = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.87'
$ws.Range('D11').Style = $style
$ws.Range('E11').Value = '  +3.16%  '
$ws.Range('E12').Value = '  -2.04%  '
$ws.Range('E13').Value = '  -2.27%  '
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '3.673.58'
$ws.Range('E15').Value = '  +0.77%  '
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('D18').Value = '63.904.65'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = '3.149.06'
$ws.Range('E19').Value = '  +0.71%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '467.49'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('E21').Value = '  +0.12%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.731'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.99'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -2.53%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.33'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +6.50%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.23'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('E28').Value = '  +8.30%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.41'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +7.61%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +0.17%  '
$style = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.69'
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  +1.95%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.110'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +1.27%  '
$ws.Range('D35').Value = '0.0₃0838'
$ws.Range('E35').Value = '  -4.58%  '
$ws.Range('E36').Value = '  +1.45%  '
$ws.Range('E37').Value = '  -2.41%  '
$ws.Range('E38').Value = '  +0.30%  '
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.26'
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  -5.45%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '463.86'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +1.04%  '
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.41'
$ws.Range('D41').Style = $style
$ws.Range('E41').Value = '  +0.73%  '
$ws.Range('E42').Value = '  +5.35%  '
$ws.Range('E43').Value = '  +5.43%  '
$ws.Range('D44').Value = '2.925.78'
$ws.Range('E44').Value = '  +0.43%  '
$ws.Range('E45').Value = '  -0.65%  '
$style = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.20'
$ws.Range('D46').Style = $style
$ws.Range('E46').Value = '  +12.66%  '
$ws.Range('E47').Value = '  -2.47%  '
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.07'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('E50').Value = '  +2.46%  '
$ws.Range('E51').Value = '  -0.81%  '
